# Apply cell-value updates to Sheet1 (cryptos list refresh).
# Values that are unambiguous text (contain non-numeric chars, e.g. two
# dots or a "%" sign) are written directly. Values that LOOK like a plain
# number (e.g. "14.30", "0.0920") are written with a leading apostrophe so
# Excel stores them as literal text instead of re-parsing/truncating them
# as a Double (which would turn "14.30" into 14.3). The cell Style is then
# reset to "Normal" so we do not leave a stray quote-prefixed number format
# behind - only the cell value changes, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.695.42"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "2.204.78"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'252.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").Value = "'0.617"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "'75.56"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -3.84%  "
$ws.Range("D10").Value = "'40.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").Value = "'0.0920"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "2.535.78"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "'14.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "2.209.63"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "'0.779"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.64%  "
$ws.Range("D18").Value = "42.638.55"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").Value = "'71.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'5.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "'2.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'228.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").Value = "'9.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.04%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'10.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").Value = "'3.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "'38.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.52%  "
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "'173.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "'20.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").Value = "'0.0829"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.63%  "
$ws.Range("D34").Value = "'5.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("D37").Value = "'0.0345"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.04%  "
$ws.Range("D38").Value = "'4.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").Value = "'12.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.59%  "
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "'2.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.86%  "
$ws.Range("D42").Value = "'59.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'5.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.69%  "
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "'102.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("D47").Value = "'8.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D48").Value = "'0.461"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.43%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("E51").Value = "  -0.87%  "
